$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 4541
$ws.Range("I101").Value = 4789.3335
$ws.Range("K101").Value = 14368.0005
$ws.Range("M101").Value = -12746.0005

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 27930.13
$ws.Range("I32").Value = 6806.684
$ws.Range("J32").Value = 128266.5
$ws.Range("K32").Value = 6806.684
$ws.Range("L32").Value = 128266.5
$ws.Range("M32").Value = -6519.684
$ws.Range("N32").Value = -128840.5
$ws.Range("H44").Value = 25324.5
$ws.Range("J44").Value = 25324.5
$ws.Range("L44").Value = 25324.5
$ws.Range("N44").Value = -26300.5
$ws.Range("H55").Value = 26431.143
$ws.Range("J55").Value = 26431.143
$ws.Range("L55").Value = 26431.143
$ws.Range("N55").Value = -27061.143
$ws.Range("H80").Value = 24081.166
$ws.Range("J80").Value = 25379.455
$ws.Range("L80").Value = 25379.455
$ws.Range("N80").Value = -27375.455
$ws.Range("H83").Value = 24081.166
$ws.Range("J83").Value = 25379.455
$ws.Range("L83").Value = 76138.36500000001
$ws.Range("N83").Value = -86122.36500000001
$ws.Range("H132").Value = 5429.091
$ws.Range("I132").Value = 5524.6665
$ws.Range("K132").Value = 16573.9995
$ws.Range("M132").Value = -14043.9995

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1325
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 1325
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 1325
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -3321
$ws.Range("H83").Value = 1325
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 1325
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 6625
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -16609
$ws.Range("H99").Value = 2608.6428
$ws.Range("I99").Value = 1376.25
$ws.Range("J99").Value = 4251.8335
$ws.Range("K99").Value = 1376.25
$ws.Range("L99").Value = 4251.8335
$ws.Range("M99").Value = 121.75
$ws.Range("N99").Value = -7247.8335
$ws.Range("H107").Value = 4703.294
$ws.Range("I107").Value = 4552.5
$ws.Range("J107").Value = 5065.2
$ws.Range("K107").Value = 4552.5
$ws.Range("L107").Value = 5065.2
$ws.Range("M107").Value = -2632.5
$ws.Range("N107").Value = -8905.200000000001

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49557.25
$ws.Range("J20").Value = 49557.25
$ws.Range("L20").Value = 49557.25
$ws.Range("N20").Value = -50029.25
$ws.Range("H30").Value = 49557.25
$ws.Range("J30").Value = 49557.25
$ws.Range("L30").Value = 49557.25
$ws.Range("N30").Value = -49739.25
$ws.Range("H41").Value = 16950
$ws.Range("J41").Value = 16950
$ws.Range("L41").Value = 16950
$ws.Range("N41").Value = -17806
$ws.Range("H51").Value = 8120
$ws.Range("J51").Value = 9326.666999999999
$ws.Range("L51").Value = 9326.666999999999
$ws.Range("N51").Value = -10798.667
$ws.Range("H60").Value = 10685.167
$ws.Range("J60").Value = 10685.167
$ws.Range("L60").Value = 10685.167
$ws.Range("N60").Value = -11707.167
$ws.Range("H61").Value = 8120
$ws.Range("J61").Value = 9326.666999999999
$ws.Range("L61").Value = 9326.666999999999
$ws.Range("N61").Value = -10022.667
$ws.Range("H109").Value = 26390
$ws.Range("J109").Value = 26390
$ws.Range("L109").Value = 26390
$ws.Range("N109").Value = -28470
$ws.Range("H123").Value = 46280
$ws.Range("J123").Value = 46280
$ws.Range("L123").Value = 46280
$ws.Range("N123").Value = -56080
$ws.Range("H128").Value = 49557.25
$ws.Range("J128").Value = 49557.25
$ws.Range("L128").Value = 49557.25
$ws.Range("N128").Value = -59517.25
$ws.Range("H132").Value = 4633.3335
$ws.Range("I132").Value = 5360
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 16080
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -13550
$ws.Range("N132").Value = -8060

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 3730
$ws.Range("I134").Value = 3529
$ws.Range("K134").Value = 10587
$ws.Range("M134").Value = -5517

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 1025.5
$ws.Range("I41").Value = 1025.5
$ws.Range("K41").Value = 1025.5
$ws.Range("M41").Value = -670.5
$ws.Range("H57").Value = 11142.875
$ws.Range("J57").Value = 14530.5
$ws.Range("L57").Value = 14530.5
$ws.Range("N57").Value = -16170.5
$ws.Range("H97").Value = 2601.6
$ws.Range("J97").Value = 2000
$ws.Range("L97").Value = 2000
$ws.Range("N97").Value = -2992
$ws.Range("H130").Value = 55588.57
$ws.Range("J130").Value = 55588.57
$ws.Range("L130").Value = 55588.57
$ws.Range("N130").Value = -65628.57000000001
$ws.Range("H132").Value = 3278.2856
$ws.Range("I132").Value = 2971.2856
$ws.Range("J132").Value = 3585.2856
$ws.Range("K132").Value = 8913.856800000001
$ws.Range("L132").Value = 10755.8568
$ws.Range("M132").Value = -6383.856800000001
$ws.Range("N132").Value = -15815.8568

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 848.7959
$ws.Range("I22").Value = 447.66666
$ws.Range("J22").Value = 1081.7097
$ws.Range("K22").Value = 447.66666
$ws.Range("L22").Value = 1081.7097
$ws.Range("M22").Value = -152.66666
$ws.Range("N22").Value = -1671.7097
$ws.Range("H27").Value = 848.7959
$ws.Range("I27").Value = 447.66666
$ws.Range("J27").Value = 1081.7097
$ws.Range("K27").Value = 447.66666
$ws.Range("L27").Value = 1081.7097
$ws.Range("M27").Value = -340.66666
$ws.Range("N27").Value = -1295.7097
$ws.Range("H46").Value = 1350.0667
$ws.Range("I46").Value = 1172.8182
$ws.Range("K46").Value = 1172.8182
$ws.Range("M46").Value = -984.8181999999999
$ws.Range("H55").Value = 581.8889
$ws.Range("I55").Value = 261.83334
$ws.Range("K55").Value = 261.83334
$ws.Range("M55").Value = -88.83334000000002
$ws.Range("H68").Value = 2452.353
$ws.Range("I68").Value = 1999.2858
$ws.Range("J68").Value = 2769.5
$ws.Range("K68").Value = 1999.2858
$ws.Range("L68").Value = 2769.5
$ws.Range("M68").Value = -1250.2858
$ws.Range("N68").Value = -4267.5
$ws.Range("H71").Value = 2452.353
$ws.Range("I71").Value = 1999.2858
$ws.Range("J71").Value = 2769.5
$ws.Range("K71").Value = 9996.429
$ws.Range("L71").Value = 13847.5
$ws.Range("M71").Value = -6252.429
$ws.Range("N71").Value = -21335.5
$ws.Range("H109").Value = 21010
$ws.Range("J109").Value = 21010
$ws.Range("L109").Value = 21010
$ws.Range("N109").Value = -23784
$ws.Range("H122").Value = 3158.7778
$ws.Range("I122").Value = 2755.9
$ws.Range("J122").Value = 3662.375
$ws.Range("K122").Value = 8267.700000000001
$ws.Range("L122").Value = 10987.125
$ws.Range("M122").Value = -5817.700000000001
$ws.Range("N122").Value = -15887.125

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3599.353
$ws.Range("J62").Value = 3645.3076
$ws.Range("L62").Value = 3645.3076
$ws.Range("N62").Value = -4893.3076
$ws.Range("H65").Value = 3599.353
$ws.Range("J65").Value = 3645.3076
$ws.Range("L65").Value = 18226.538
$ws.Range("N65").Value = -24466.538
$ws.Range("H109").Value = 17838
$ws.Range("J109").Value = 17838
$ws.Range("L109").Value = 17838
$ws.Range("N109").Value = -20612
$ws.Range("H122").Value = 1985.7667
$ws.Range("I122").Value = 1669.4073
$ws.Range("J122").Value = 4833
$ws.Range("K122").Value = 5008.2219
$ws.Range("L122").Value = 14499
$ws.Range("M122").Value = -2558.2219
$ws.Range("N122").Value = -19399
$ws.Range("H136").Value = 1153.1428
$ws.Range("I136").Value = 918.129
$ws.Range("J136").Value = 1815.4546
$ws.Range("K136").Value = 2754.387
$ws.Range("L136").Value = 5446.3638
$ws.Range("M136").Value = -204.3870000000002
$ws.Range("N136").Value = -10546.3638

Write-Output "Applied all profit sheet updates"